$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces literal text entry (mirrors typing '<value> into Excel),
# preventing numeric-looking strings (e.g. "581.58", "1.00") from being coerced to numbers.

$ws.Range("D2").Value = "'63.058.60"
$ws.Range("E2").Value = "'  -0.74%  "
$ws.Range("D3").Value = "'2.551.44"
$ws.Range("E3").Value = "'  +0.28%  "
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("D5").Value = "'581.58"
$ws.Range("E5").Value = "'  +2.22%  "
$ws.Range("D6").Value = "'147.06"
$ws.Range("E6").Value = "'  -2.31%  "
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E8").Value = "'  -0.30%  "
$ws.Range("E9").Value = "'  +0.07%  "
$ws.Range("E10").Value = "'  -2.38%  "
$ws.Range("E11").Value = "'  -0.01%  "
$ws.Range("E12").Value = "'  -1.30%  "
$ws.Range("D13").Value = "'27.54"
$ws.Range("E13").Value = "'  -2.48%  "
$ws.Range("D14").Value = "'3.007.01"
$ws.Range("E14").Value = "'  +0.18%  "
$ws.Range("D15").Value = "'62.969.97"
$ws.Range("E15").Value = "'  -0.74%  "
$ws.Range("E16").Value = "'  +0.16%  "
$ws.Range("D17").Value = "'2.548.00"
$ws.Range("E17").Value = "'  +0.11%  "
$ws.Range("D18").Value = "'11.34"
$ws.Range("E18").Value = "'  -2.27%  "
$ws.Range("D19").Value = "'338.79"
$ws.Range("E19").Value = "'  -0.19%  "
$ws.Range("E20").Value = "'  -0.64%  "
$ws.Range("D21").Value = "'6.75"
$ws.Range("E21").Value = "'  -0.80%  "
$ws.Range("E22").Value = "'  -0.09%  "
$ws.Range("D23").Value = "'65.51"
$ws.Range("E23").Value = "'  -0.74%  "
$ws.Range("D24").Value = "'2.677.79"
$ws.Range("E24").Value = "'  +0.27%  "
$ws.Range("E25").Value = "'  -0.19%  "
$ws.Range("D26").Value = "'1.61"
$ws.Range("E26").Value = "'  +1.02%  "
$ws.Range("E27").Value = "'  -4.70%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "'  -0.04%  "
$ws.Range("D29").Value = "'8.36"
$ws.Range("E29").Value = "'  -1.26%  "
$ws.Range("D30").Value = "'7.70"
$ws.Range("E30").Value = "'  +7.71%  "
$ws.Range("D31").Value = "'1.99"
$ws.Range("E31").Value = "'  +5.54%  "
$ws.Range("E32").Value = "'  -0.89%  "
$ws.Range("D33").Value = "'178.06"
$ws.Range("D34").Value = "'421.77"
$ws.Range("E34").Value = "'  +0.39%  "
$ws.Range("E35").Value = "'  -1.90%  "
$ws.Range("D36").Value = "'0.400"
$ws.Range("E36").Value = "'  -1.42%  "
$ws.Range("D37").Value = "'19.09"
$ws.Range("E37").Value = "'  +0.21%  "
$ws.Range("E39").Value = "'  -1.37%  "
$ws.Range("E40").Value = "'  -2.06%  "
$ws.Range("E41").Value = "'  +0.00%  "
$ws.Range("D42").Value = "'39.78"
$ws.Range("E42").Value = "'  +0.62%  "
$ws.Range("D43").Value = "'150.80"
$ws.Range("E43").Value = "'  -1.92%  "
$ws.Range("D44").Value = "'3.77"
$ws.Range("E44").Value = "'  -0.25%  "
$ws.Range("D45").Value = "'20.77"
$ws.Range("E45").Value = "'  -0.28%  "
$ws.Range("E46").Value = "'  +2.07%  "
$ws.Range("E47").Value = "'  -1.30%  "
$ws.Range("E48").Value = "'  +0.36%  "
$ws.Range("E49").Value = "'  +0.04%  "
$ws.Range("D50").Value = "'18.30"
$ws.Range("E50").Value = "'  -1.84%  "
$ws.Range("E51").Value = "'  -5.96%  "
